$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("MNT-USD", "IMX-USD", "TAO-USD", "GRT-USD", "PEPE-USD")

$startRow = 489
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
